$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "IT - OPEN ACTION ITEM LOG"
$ws.Range("A2").Value = "Project: IT RAID Log"
$ws.Range("E12").Value = "Ethics Committee"
$ws.Range("E13").Value = "Compliance Officers"

$wb.Save()
